$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.989.04"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +6.40%  "

# Row 3
$ws.Range("D3").Value = "'1.884.23"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +5.63%  "

# Row 4
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'248.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.02%  "

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("D7").Value = "'0.4982"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.34%  "

# Row 8
$ws.Range("D8").Value = "'45.87"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +9.29%  "

# Row 9
$ws.Range("D9").Value = "'0.2858"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.98%  "

# Row 10
$ws.Range("D10").Value = "'0.06544"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.65%  "

# Row 11
$ws.Range("D11").Value = "'1.881.20"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.52%  "

# Row 12
$ws.Range("E12").Value = "  +3.86%  "

# Row 13
$ws.Range("D13").Value = "'0.07213"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.74%  "

# Row 14
$ws.Range("D14").Value = "'0.6640"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.16%  "

# Row 15
$ws.Range("D15").Value = "'85.12"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.48%  "

# Row 16
$ws.Range("D16").Value = "'4.798"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.57%  "

# Row 17
$ws.Range("D17").Value = "'30.002.68"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.62%  "

# Row 18
$ws.Range("D18").Value = "'0.9995"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.01%  "

# Row 19
$ws.Range("D19").Value = "'12.88"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.07%  "

# Row 20
$ws.Range("D20").Value = "'0.000007497"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.97%  "

# Row 21
$ws.Range("D21").Value = "'0.9991"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.06%  "

# Row 22
$ws.Range("D22").Value = "'2.124.56"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.81%  "

# Row 23
$ws.Range("D23").Value = "'4.752"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.24%  "

# Row 24
$ws.Range("D24").Value = "'5.541"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.96%  "

# Row 25
$ws.Range("D25").Value = "'9.010"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.31%  "

# Row 26
$ws.Range("D26").Value = "'145.08"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.59%  "

# Row 27
$ws.Range("D27").Value = "'134.67"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +23.60%  "

# Row 28
$ws.Range("E28").Value = "  +5.78%  "

# Row 29
$ws.Range("E29").Value = "  +5.35%  "

# Row 30
$ws.Range("D30").Value = "'1.377"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.70%  "

# Row 31
$ws.Range("D31").Value = "'4.171"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.10%  "

# Row 32
$ws.Range("D32").Value = "'0.08616"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.24%  "

# Row 33
$ws.Range("D33").Value = "'3.870"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.44%  "

# Row 34
$ws.Range("D34").Value = "'0.05115"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.57%  "

# Row 35
$ws.Range("D35").Value = "'1.130"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.57%  "

# Row 36
$ws.Range("D36").Value = "'0.6870"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.63%  "

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.709"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.69%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.312"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +13.14%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.753"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.45%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.9581"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.66%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.01631"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.32%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.074"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.24%  "

# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'104.29"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.61%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.0000"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4217"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.96%  "

# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.431"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.54%  "

# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1253"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.41%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05635"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.71%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'32.37"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.84%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.245"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.42%  "

# Row 51
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3717"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.18%  "
